$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update feature names (reordering rows 5-9) and importance values (rows 2-10)
$ws.Range("A2").Value = "RSI"
$ws.Range("B2").Value = 0.4605642412342682

$ws.Range("A3").Value = "MACD"
$ws.Range("B3").Value = 0.2449307532255113

$ws.Range("A4").Value = "Signal_line"
$ws.Range("B4").Value = 0.09795370241457584

$ws.Range("A5").Value = "close_long"
$ws.Range("B5").Value = 0.04387137303273737

$ws.Range("A6").Value = "close_short"
$ws.Range("B6").Value = 0.03784587917419765

$ws.Range("A7").Value = "VIX_short"
$ws.Range("B7").Value = 0.03443029589432096

$ws.Range("A8").Value = "VIX"
$ws.Range("B8").Value = 0.02745765412396424

$ws.Range("A9").Value = "VIX_long"
$ws.Range("B9").Value = 0.02647626597958512

$ws.Range("A10").Value = "DJI"
$ws.Range("B10").Value = 0.02646983492083943
